$d = $word.ActiveDocument

# Each (old -> new) multiplication expression replacement, derived from the diff.
# All old values are unique across the document, so a straightforward
# Find/Replace (MatchCase, MatchWholeWord off, with trailing "=" to anchor)
# for each pair is unambiguous and safe to run sequentially.

$range = $d.Content
$range.Find.Execute("99×88=", $true, $true, $false, $false, $false, $true, 1, $false, "37×14=", 2) | Out-Null

$range = $d.Content
$range.Find.Execute("15×45=", $true, $true, $false, $false, $false, $true, 1, $false, "89×52=", 2) | Out-Null

$range = $d.Content
$range.Find.Execute("89×27=", $true, $true, $false, $false, $false, $true, 1, $false, "71×97=", 2) | Out-Null

$range = $d.Content
$range.Find.Execute("13×17=", $true, $true, $false, $false, $false, $true, 1, $false, "98×52=", 2) | Out-Null

$range = $d.Content
$range.Find.Execute("91×52=", $true, $true, $false, $false, $false, $true, 1, $false, "50×52=", 2) | Out-Null

$range = $d.Content
$range.Find.Execute("66×74=", $true, $true, $false, $false, $false, $true, 1, $false, "98×93=", 2) | Out-Null

$range = $d.Content
$range.Find.Execute("99×83=", $true, $true, $false, $false, $false, $true, 1, $false, "25×43=", 2) | Out-Null

$range = $d.Content
$range.Find.Execute("96×88=", $true, $true, $false, $false, $false, $true, 1, $false, "28×20=", 2) | Out-Null

$range = $d.Content
$range.Find.Execute("18×68=", $true, $true, $false, $false, $false, $true, 1, $false, "98×84=", 2) | Out-Null

$range = $d.Content
$range.Find.Execute("34×47=", $true, $true, $false, $false, $false, $true, 1, $false, "35×72=", 2) | Out-Null

$range = $d.Content
$range.Find.Execute("93×64=", $true, $true, $false, $false, $false, $true, 1, $false, "48×51=", 2) | Out-Null

$range = $d.Content
$range.Find.Execute("31×18=", $true, $true, $false, $false, $false, $true, 1, $false, "98×95=", 2) | Out-Null

$range = $d.Content
$range.Find.Execute("84×26=", $true, $true, $false, $false, $false, $true, 1, $false, "52×41=", 2) | Out-Null

$range = $d.Content
$range.Find.Execute("78×34=", $true, $true, $false, $false, $false, $true, 1, $false, "33×45=", 2) | Out-Null

$range = $d.Content
$range.Find.Execute("83×40=", $true, $true, $false, $false, $false, $true, 1, $false, "98×63=", 2) | Out-Null

$range = $d.Content
$range.Find.Execute("18×51=", $true, $true, $false, $false, $false, $true, 1, $false, "90×63=", 2) | Out-Null

$range = $d.Content
$range.Find.Execute("88×67=", $true, $true, $false, $false, $false, $true, 1, $false, "94×34=", 2) | Out-Null

$range = $d.Content
$range.Find.Execute("41×26=", $true, $true, $false, $false, $false, $true, 1, $false, "54×99=", 2) | Out-Null

$range = $d.Content
$range.Find.Execute("29×45=", $true, $true, $false, $false, $false, $true, 1, $false, "53×50=", 2) | Out-Null

$range = $d.Content
$range.Find.Execute("28×45=", $true, $true, $false, $false, $false, $true, 1, $false, "87×67=", 2) | Out-Null

$range = $d.Content
$range.Find.Execute("35×45=", $true, $true, $false, $false, $false, $true, 1, $false, "78×57=", 2) | Out-Null

$range = $d.Content
$range.Find.Execute("73×19=", $true, $true, $false, $false, $false, $true, 1, $false, "91×54=", 2) | Out-Null

$range = $d.Content
$range.Find.Execute("29×13=", $true, $true, $false, $false, $false, $true, 1, $false, "65×78=", 2) | Out-Null

$range = $d.Content
$range.Find.Execute("40×78=", $true, $true, $false, $false, $false, $true, 1, $false, "80×51=", 2) | Out-Null

$range = $d.Content
$range.Find.Execute("37×66=", $true, $true, $false, $false, $false, $true, 1, $false, "76×55=", 2) | Out-Null

$d.Save()